# Apply numeric updates to the Behemoth_Profits workbook (Leve profit recompute)
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1069125.2
$ws.Range("I17").Value = 9259353
$ws.Range("J17").Value = 834.6957
$ws.Range("K17").Value = 27778059
$ws.Range("L17").Value = 2504.0871
$ws.Range("M17").Value = -27777891
$ws.Range("N17").Value = -2840.0871
$ws.Range("H55").Value = 3519.6667
$ws.Range("J55").Value = 19999
$ws.Range("L55").Value = 19999
$ws.Range("N55").Value = -20427
$ws.Range("H76").Value = 5404.1333
$ws.Range("I76").Value = 3673.6667
$ws.Range("J76").Value = 7999.8335
$ws.Range("K76").Value = 3673.6667
$ws.Range("L76").Value = 7999.8335
$ws.Range("M76").Value = -3358.6667
$ws.Range("N76").Value = -8629.833500000001
$ws.Range("H79").Value = 5404.1333
$ws.Range("I79").Value = 3673.6667
$ws.Range("J79").Value = 7999.8335
$ws.Range("K79").Value = 3673.6667
$ws.Range("L79").Value = 7999.8335
$ws.Range("M79").Value = -2581.6667
$ws.Range("N79").Value = -10183.8335
$ws.Range("H100").Value = 2778.375
$ws.Range("I100").Value = 1327.2
$ws.Range("K100").Value = 1327.2
$ws.Range("M100").Value = -786.2
$ws.Range("H113").Value = 100003400
$ws.Range("I113").Value = 33337000
$ws.Range("K113").Value = 33337000
$ws.Range("M113").Value = -33333746
$ws.Range("H132").Value = 1076.4445
$ws.Range("I132").Value = 1010.75
$ws.Range("J132").Value = 1602
$ws.Range("K132").Value = 3032.25
$ws.Range("L132").Value = 4806
$ws.Range("M132").Value = -502.25
$ws.Range("N132").Value = -9866
$ws.Range("H133").Value = 60795.832
$ws.Range("J133").Value = 60795.832
$ws.Range("L133").Value = 60795.832
$ws.Range("N133").Value = -70915.83199999999

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9435619
$ws.Range("J32").Value = 2430.1667
$ws.Range("L32").Value = 2430.1667
$ws.Range("N32").Value = -3004.1667
$ws.Range("H45").Value = 2349.625
$ws.Range("I45").Value = 2002
$ws.Range("K45").Value = 2002
$ws.Range("M45").Value = -1625
$ws.Range("H61").Value = 31318282
$ws.Range("J61").Value = 133491
$ws.Range("L61").Value = 133491
$ws.Range("N61").Value = -133915
$ws.Range("H74").Value = 10425766
$ws.Range("I74").Value = 20835620
$ws.Range("K74").Value = 20835620
$ws.Range("M74").Value = -20834746
$ws.Range("H77").Value = 10425766
$ws.Range("I77").Value = 20835620
$ws.Range("K77").Value = 104178100
$ws.Range("M77").Value = -104173732
$ws.Range("H88").Value = 1224.0625
$ws.Range("I88").Value = 783.625
$ws.Range("K88").Value = 783.625
$ws.Range("M88").Value = -377.625
$ws.Range("H91").Value = 1224.0625
$ws.Range("I91").Value = 783.625
$ws.Range("K91").Value = 783.625
$ws.Range("M91").Value = 620.375
$ws.Range("H102").Value = 3228.2856
$ws.Range("I102").Value = 3228.2856
$ws.Range("K102").Value = 3228.2856
$ws.Range("M102").Value = -1606.2856
$ws.Range("H106").Value = 28353.5
$ws.Range("J106").Value = 28353.5
$ws.Range("L106").Value = 28353.5
$ws.Range("N106").Value = -30877.5
$ws.Range("H110").Value = 2356.353
$ws.Range("I110").Value = 2228.75
$ws.Range("J110").Value = 2662.6
$ws.Range("K110").Value = 2228.75
$ws.Range("L110").Value = 2662.6
$ws.Range("M110").Value = -183.75
$ws.Range("N110").Value = -6752.6
$ws.Range("H132").Value = 2760.5
$ws.Range("I132").Value = 2760.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8281.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5751.5
$ws.Range("N132").Value = $null
$ws.Range("H136").Value = 31318282
$ws.Range("J136").Value = 133491
$ws.Range("L136").Value = 400473
$ws.Range("N136").Value = -405573

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1099.9474
$ws.Range("I16").Value = 760.5333000000001
$ws.Range("K16").Value = 760.5333000000001
$ws.Range("M16").Value = -473.5333000000001
$ws.Range("H113").Value = 1099.9474
$ws.Range("I113").Value = 760.5333000000001
$ws.Range("K113").Value = 760.5333000000001
$ws.Range("M113").Value = 1409.4667
$ws.Range("H132").Value = 2577.7856
$ws.Range("I132").Value = 2632.1482
$ws.Range("J132").Value = 1110
$ws.Range("K132").Value = 7896.444600000001
$ws.Range("L132").Value = 3330
$ws.Range("M132").Value = -5366.444600000001
$ws.Range("N132").Value = -8390

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1040.7222
$ws.Range("I140").Value = 1040.7222
$ws.Range("K140").Value = 3122.1666
$ws.Range("M140").Value = 2057.8334

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 96432.28999999999
$ws.Range("J62").Value = 98989.8
$ws.Range("L62").Value = 98989.8
$ws.Range("N62").Value = -100361.8
$ws.Range("H65").Value = 96432.28999999999
$ws.Range("J65").Value = 98989.8
$ws.Range("L65").Value = 296969.4
$ws.Range("N65").Value = -303833.4
$ws.Range("H132").Value = 30304782
$ws.Range("I132").Value = 30304782
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 90914346
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -90911816
$ws.Range("N132").Value = $null

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2995.9412
$ws.Range("I46").Value = 2433
$ws.Range("K46").Value = 2433
$ws.Range("M46").Value = -2245
$ws.Range("H61").Value = 6798.778
$ws.Range("I61").Value = 1998.3334
$ws.Range("K61").Value = 1998.3334
$ws.Range("M61").Value = -1796.3334
$ws.Range("H113").Value = 6798.778
$ws.Range("I113").Value = 1998.3334
$ws.Range("K113").Value = 1998.3334
$ws.Range("M113").Value = 171.6666
$ws.Range("H132").Value = 253814.88
$ws.Range("I132").Value = 254628.5
$ws.Range("K132").Value = 763885.5
$ws.Range("M132").Value = -761355.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 39999
$ws.Range("I38").Value = 39999
$ws.Range("K38").Value = 39999
$ws.Range("M38").Value = -39526
$ws.Range("H49").Value = 59999
$ws.Range("I49").Value = 59999
$ws.Range("K49").Value = 59999
$ws.Range("M49").Value = -59769
$ws.Range("H126").Value = 5798.6553
$ws.Range("I126").Value = 5885.4346
$ws.Range("K126").Value = 17656.3038
$ws.Range("M126").Value = -15186.3038
